$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 6231.544
$ws.Range("I15").Value = 6231.544
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 18694.632
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -18525.632

# Row 70
$ws.Range("H70").Value = 1464.5834
$ws.Range("I70").Value = 1333.3334
$ws.Range("J70").Value = 1508.3334
$ws.Range("K70").Value = 4000.0002
$ws.Range("L70").Value = 4525.0002
$ws.Range("M70").Value = -3730.0002
$ws.Range("N70").Value = -5065.0002

# Row 73
$ws.Range("H73").Value = 1464.5834
$ws.Range("I73").Value = 1333.3334
$ws.Range("J73").Value = 1508.3334
$ws.Range("K73").Value = 4000.0002
$ws.Range("L73").Value = 4525.0002
$ws.Range("M73").Value = -3064.0002
$ws.Range("N73").Value = -6397.0002

# Row 80
$ws.Range("H80").Value = 20835060
$ws.Range("I80").Value = 55556550
$ws.Range("J80").Value = 2164.8
$ws.Range("K80").Value = 166669650
$ws.Range("L80").Value = 6494.400000000001
$ws.Range("M80").Value = -166668652
$ws.Range("N80").Value = -8490.400000000001

# Row 83
$ws.Range("H83").Value = 20835060
$ws.Range("I83").Value = 55556550
$ws.Range("J83").Value = 2164.8
$ws.Range("K83").Value = 500008950
$ws.Range("L83").Value = 19483.2
$ws.Range("M83").Value = -500003958
$ws.Range("N83").Value = -29467.2

# Row 112
$ws.Range("H112").Value = 473421.3
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 491967.25
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 1475901.75
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -1478117.75

# Row 125
$ws.Range("H125").Value = 1188.3636
$ws.Range("I125").Value = 1137
$ws.Range("J125").Value = 1217.7142
$ws.Range("K125").Value = 10233
$ws.Range("L125").Value = 10959.4278
$ws.Range("M125").Value = -7773
$ws.Range("N125").Value = -15879.4278

# Row 127
$ws.Range("H127").Value = 1328
$ws.Range("I127").Value = 1041
$ws.Range("J127").Value = 1650.875
$ws.Range("K127").Value = 3123
$ws.Range("L127").Value = 4952.625
$ws.Range("M127").Value = 1837
$ws.Range("N127").Value = -14872.625

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 11388.833
$ws.Range("I25").Value = 333
$ws.Range("J25").Value = 13600
$ws.Range("K25").Value = 333
$ws.Range("L25").Value = 13600
$ws.Range("M25").Value = 69
$ws.Range("N25").Value = -14404

# Row 30
$ws.Range("H30").Value = 19995.5
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 19995.5
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 19995.5
$ws.Range("N30").Value = -20295.5
$ws.Range("M30").ClearContents()

# Row 32
$ws.Range("H32").Value = 9482.358
$ws.Range("I32").Value = 8519.816000000001
$ws.Range("J32").Value = 17732.715
$ws.Range("K32").Value = 8519.816000000001
$ws.Range("L32").Value = 17732.715
$ws.Range("M32").Value = -8232.816000000001
$ws.Range("N32").Value = -18306.715

# Row 110
$ws.Range("H110").Value = 1487.2122
$ws.Range("I110").Value = 1534.7742
$ws.Range("J110").Value = 750
$ws.Range("K110").Value = 1534.7742
$ws.Range("L110").Value = 750
$ws.Range("M110").Value = 510.2257999999999
$ws.Range("N110").Value = -4840

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5380745.5
$ws.Range("I31").Value = 1794.129
$ws.Range("J31").Value = 10759697
$ws.Range("K31").Value = 1794.129
$ws.Range("L31").Value = 10759697
$ws.Range("M31").Value = -1499.129
$ws.Range("N31").Value = -10760287

# Row 34
$ws.Range("H34").Value = 5380745.5
$ws.Range("I34").Value = 1794.129
$ws.Range("J34").Value = 10759697
$ws.Range("K34").Value = 1794.129
$ws.Range("L34").Value = 10759697
$ws.Range("M34").Value = -1592.129
$ws.Range("N34").Value = -10760101

# Row 99
$ws.Range("H99").Value = 3017.1365
$ws.Range("I99").Value = 2815.375
$ws.Range("J99").Value = 3132.4285
$ws.Range("K99").Value = 2815.375
$ws.Range("L99").Value = 3132.4285
$ws.Range("M99").Value = -1317.375
$ws.Range("N99").Value = -6128.4285

# Row 100
$ws.Range("H100").Value = 25000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 25000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164

# Row 126
$ws.Range("H126").Value = 3017.1365
$ws.Range("I126").Value = 2815.375
$ws.Range("J126").Value = 3132.4285
$ws.Range("K126").Value = 8446.125
$ws.Range("L126").Value = 9397.2855
$ws.Range("M126").Value = -5976.125
$ws.Range("N126").Value = -14337.2855

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5232.909
$ws.Range("I5").Value = 7214.933
$ws.Range("J5").Value = 985.7143
$ws.Range("K5").Value = 21644.799
$ws.Range("L5").Value = 2957.1429
$ws.Range("M5").Value = -21532.799
$ws.Range("N5").Value = -3181.1429

# Row 75
$ws.Range("H75").Value = 2540
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 2711.111
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 8133.333
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -10129.333

# Row 78
$ws.Range("H78").Value = 2540
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 2711.111
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 24399.999
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -34383.999

# Row 122
$ws.Range("H122").Value = 6714.35
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 7030.8945
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 63278.0505
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -68178.05050000001

# Row 132
$ws.Range("H132").Value = 4123.3887
$ws.Range("I132").Value = 1456.25
$ws.Range("J132").Value = 6257.1
$ws.Range("K132").Value = 13106.25
$ws.Range("L132").Value = 56313.9
$ws.Range("M132").Value = -10576.25
$ws.Range("N132").Value = -61373.9

# Row 135
$ws.Range("H135").Value = 5232.909
$ws.Range("I135").Value = 7214.933
$ws.Range("J135").Value = 985.7143
$ws.Range("K135").Value = 64934.397
$ws.Range("L135").Value = 8871.4287
$ws.Range("M135").Value = -62399.397
$ws.Range("N135").Value = -13941.4287

# Row 138
$ws.Range("H138").Value = 2110.9092
$ws.Range("I138").Value = 1522
$ws.Range("J138").Value = 8000
$ws.Range("K138").Value = 4566
$ws.Range("L138").Value = 24000
$ws.Range("M138").Value = 574
$ws.Range("N138").Value = -34280

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 34447690
$ws.Range("I24").Value = 103333336
$ws.Range("J24").Value = 4864.5
$ws.Range("K24").Value = 103333336
$ws.Range("L24").Value = 4864.5
$ws.Range("M24").Value = -103333163
$ws.Range("N24").Value = -5210.5

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5766.9165
$ws.Range("I7").Value = 2034.6666
$ws.Range("J7").Value = 9499.166999999999
$ws.Range("K7").Value = 2034.6666
$ws.Range("L7").Value = 9499.166999999999
$ws.Range("M7").Value = -1922.6666
$ws.Range("N7").Value = -9723.166999999999

# Row 122
$ws.Range("H122").Value = 2211.1
$ws.Range("I122").Value = 2234.5557

# Row 126
$ws.Range("H126").Value = 5766.9165
$ws.Range("I126").Value = 2034.6666
$ws.Range("J126").Value = 9499.166999999999
$ws.Range("K126").Value = 6103.9998
$ws.Range("L126").Value = 28497.501
$ws.Range("M126").Value = -3633.9998
$ws.Range("N126").Value = -33437.501

# Row 132
$ws.Range("H132").Value = 4695.3477
$ws.Range("I132").Value = 4333.4165
$ws.Range("J132").Value = 5090.1816
$ws.Range("K132").Value = 13000.2495
$ws.Range("L132").Value = 15270.5448
$ws.Range("M132").Value = -10470.2495
$ws.Range("N132").Value = -20330.5448
